$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "58.166.73"
$r.Style = "Normal"
$ws.Range("E2").Value = "  +0.53%  "
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "3.141.92"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.04%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "534.56"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "138.64"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -0.03%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "3.140.42"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.467"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +4.33%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "7.31"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("E11").Value = "  -0.47%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.415"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +4.66%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "3.674.98"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("E14").Value = "  +1.59%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "25.69"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  +0.00%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "58.280.91"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +0.49%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "3.146.05"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +0.37%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "6.05"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +0.60%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "12.72"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -0.58%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "8.18"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +2.69%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "360.37"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +1.80%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +0.47%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "69.11"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "0.507"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("E27").Value = "  -0.05%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "0.0₃0881"
$r.Style = "Normal"
$ws.Range("E28").Value = "  -3.72%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "7.33"
$r.Style = "Normal"
$ws.Range("E29").Value = "  -2.18%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "6.18"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  +1.61%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "5.04"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("E34").Value = "  -2.98%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "159.46"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  -1.45%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "26.17"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("E38").Value = "  +0.39%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "1.70"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +4.90%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.0671"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -0.01%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "2.509.16"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +8.49%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.703"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "4.00"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -4.18%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "37.44"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +2.43%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "3.190.53"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +0.45%  "
$r = $ws.Range("B46")
$r.NumberFormat = "@"
$r.Value = "FirstDigitalUSD"
$r.Style = "Normal"
$r = $ws.Range("C46")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$r.Style = "Normal"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "
$r = $ws.Range("B47")
$r.NumberFormat = "@"
$r.Value = "VeChain"
$r.Style = "Normal"
$r = $ws.Range("C47")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$r.Style = "Normal"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.0269"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -0.88%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.990"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("E49").Value = "  +0.41%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "19.86"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("E51").Value = "  -4.17%  "
